$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(9, 8).Value = 7786.4287  # H9: was 5884.5264
$ws.Cells.Item(9, 9).Value = 8000.769  # I9: was 5933.6665
$ws.Cells.Item(9, 11).Value = 8000.769  # K9: was 5933.6665
$ws.Cells.Item(9, 13).Value = -7831.769  # M9: was -5764.6665
$ws.Cells.Item(15, 8).Value = 738.94916  # H15: was 771.7166999999999
$ws.Cells.Item(15, 9).Value = 738.94916  # I15: was 771.7166999999999
$ws.Cells.Item(15, 11).Value = 2216.84748  # K15: was 2315.1501
$ws.Cells.Item(15, 13).Value = -2047.84748  # M15: was -2146.1501
$ws.Cells.Item(70, 8).Value = 2180.1428  # H70: was 2302.4443
$ws.Cells.Item(70, 9).Value = 2350.5  # I70: was 2400.6667
$ws.Cells.Item(70, 10).Value = 1953  # J70: was 2253.3333
$ws.Cells.Item(70, 11).Value = 7051.5  # K70: was 7202.000100000001
$ws.Cells.Item(70, 12).Value = 5859  # L70: was 6759.999899999999
$ws.Cells.Item(70, 13).Value = -6781.5  # M70: was -6932.000100000001
$ws.Cells.Item(70, 14).Value = -6399  # N70: was -7299.999899999999
$ws.Cells.Item(73, 8).Value = 2180.1428  # H73: was 2302.4443
$ws.Cells.Item(73, 9).Value = 2350.5  # I73: was 2400.6667
$ws.Cells.Item(73, 10).Value = 1953  # J73: was 2253.3333
$ws.Cells.Item(73, 11).Value = 7051.5  # K73: was 7202.000100000001
$ws.Cells.Item(73, 12).Value = 5859  # L73: was 6759.999899999999
$ws.Cells.Item(73, 13).Value = -6115.5  # M73: was -6266.000100000001
$ws.Cells.Item(73, 14).Value = -7731  # N73: was -8631.999899999999
$ws.Cells.Item(76, 8).Value = 11481  # H76: was 12125.77
$ws.Cells.Item(76, 9).Value = 20664  # I76: was 20814.166
$ws.Cells.Item(76, 10).Value = 4593.75  # J76: was 4678.5713
$ws.Cells.Item(76, 11).Value = 20664  # K76: was 20814.166
$ws.Cells.Item(76, 12).Value = 4593.75  # L76: was 4678.5713
$ws.Cells.Item(76, 13).Value = -20349  # M76: was -20499.166
$ws.Cells.Item(76, 14).Value = -5223.75  # N76: was -5308.5713
$ws.Cells.Item(79, 8).Value = 11481  # H79: was 12125.77
$ws.Cells.Item(79, 9).Value = 20664  # I79: was 20814.166
$ws.Cells.Item(79, 10).Value = 4593.75  # J79: was 4678.5713
$ws.Cells.Item(79, 11).Value = 20664  # K79: was 20814.166
$ws.Cells.Item(79, 12).Value = 4593.75  # L79: was 4678.5713
$ws.Cells.Item(79, 13).Value = -19572  # M79: was -19722.166
$ws.Cells.Item(79, 14).Value = -6777.75  # N79: was -6862.5713
$ws.Cells.Item(98, 8).Value = 1348.7  # H98: was 1140.7778
$ws.Cells.Item(98, 9).Value = 1213.1428  # I98: was 1064.625
$ws.Cells.Item(98, 10).Value = 1665  # J98: was 1750
$ws.Cells.Item(98, 11).Value = 1213.1428  # K98: was 1064.625
$ws.Cells.Item(98, 12).Value = 1665  # L98: was 1750
$ws.Cells.Item(98, 13).Value = 284.8571999999999  # M98: was 433.375
$ws.Cells.Item(98, 14).Value = -4661  # N98: was -4746
$ws.Cells.Item(106, 8).Value = 1798.4  # H106: was 2224.3333
$ws.Cells.Item(106, 9).Value = 1954.5  # I106: was 2749.5
$ws.Cells.Item(106, 11).Value = 1954.5  # K106: was 2749.5
$ws.Cells.Item(106, 13).Value = -1323.5  # M106: was -2118.5
$ws.Cells.Item(116, 8).Value = 7417.9  # H116: was 7668
$ws.Cells.Item(116, 9).Value = 6502.6665  # I116: was 6753.25
$ws.Cells.Item(116, 10).Value = 8790.75  # J116: was 8887.666999999999
$ws.Cells.Item(116, 11).Value = 6502.6665  # K116: was 6753.25
$ws.Cells.Item(116, 12).Value = 8790.75  # L116: was 8887.666999999999
$ws.Cells.Item(116, 13).Value = -3060.6665  # M116: was -3311.25
$ws.Cells.Item(116, 14).Value = -15674.75  # N116: was -15771.667
$ws.Cells.Item(122, 8).Value = 1348.7  # H122: was 1140.7778
$ws.Cells.Item(122, 9).Value = 1213.1428  # I122: was 1064.625
$ws.Cells.Item(122, 10).Value = 1665  # J122: was 1750
$ws.Cells.Item(122, 11).Value = 3639.4284  # K122: was 3193.875
$ws.Cells.Item(122, 12).Value = 4995  # L122: was 5250
$ws.Cells.Item(122, 13).Value = -1189.4284  # M122: was -743.875
$ws.Cells.Item(122, 14).Value = -9895  # N122: was -10150
$ws.Cells.Item(132, 8).Value = 2617.3076  # H132: was 2683.28
$ws.Cells.Item(132, 9).Value = 2334.9583  # I132: was 2394.3914
$ws.Cells.Item(132, 11).Value = 7004.874899999999  # K132: was 7183.174199999999
$ws.Cells.Item(132, 13).Value = -4474.874899999999  # M132: was -4653.174199999999
$ws.Cells.Item(137, 8).Value = 4617.625  # H137: was 4422.4116
$ws.Cells.Item(137, 9).Value = 5427  # I137: was 5079
$ws.Cells.Item(137, 10).Value = 3808.25  # J137: was 3683.75
$ws.Cells.Item(137, 11).Value = 16281  # K137: was 15237
$ws.Cells.Item(137, 12).Value = 11424.75  # L137: was 11051.25
$ws.Cells.Item(137, 13).Value = -13731  # M137: was -12687
$ws.Cells.Item(137, 14).Value = -16524.75  # N137: was -16151.25

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 6887.3076  # H32: was 7219.324
$ws.Cells.Item(32, 9).Value = 6071.222  # I32: was 6384.5293
$ws.Cells.Item(32, 11).Value = 6071.222  # K32: was 6384.5293
$ws.Cells.Item(32, 13).Value = -5784.222  # M32: was -6097.5293
$ws.Cells.Item(122, 8).Value = 2553.5454  # H122: was 3636.4375
$ws.Cells.Item(122, 9).Value = 2009.8889  # I122: was 2433.1667
$ws.Cells.Item(122, 10).Value = 5000  # J122: was 4358.4
$ws.Cells.Item(122, 11).Value = 6029.6667  # K122: was 7299.500100000001
$ws.Cells.Item(122, 12).Value = 15000  # L122: was 13075.2
$ws.Cells.Item(122, 13).Value = -3579.6667  # M122: was -4849.500100000001
$ws.Cells.Item(122, 14).Value = -19900  # N122: was -17975.2
$ws.Cells.Item(132, 8).Value = 31302328  # H132: was 32312074
$ws.Cells.Item(132, 9).Value = 12312  # I132: was 12350.429
$ws.Cells.Item(132, 10).Value = 91037816  # J132: was 100141496
$ws.Cells.Item(132, 11).Value = 36936  # K132: was 37051.287
$ws.Cells.Item(132, 12).Value = 273113448  # L132: was 300424488
$ws.Cells.Item(132, 13).Value = -34406  # M132: was -34521.287
$ws.Cells.Item(132, 14).Value = -273118508  # N132: was -300429548

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(7, 8).Value = 188.10715  # H7: was 227.70833
$ws.Cells.Item(7, 9).Value = 49.357143  # I7: was 67.23077000000001
$ws.Cells.Item(7, 10).Value = 326.85715  # J7: was 417.36365
$ws.Cells.Item(7, 11).Value = 49.357143  # K7: was 67.23077000000001
$ws.Cells.Item(7, 12).Value = 326.85715  # L7: was 417.36365
$ws.Cells.Item(7, 13).Value = 63.642857  # M7: was 45.76922999999999
$ws.Cells.Item(7, 14).Value = -552.85715  # N7: was -643.36365
$ws.Cells.Item(31, 8).Value = 4342.591  # H31: was 4466.143
$ws.Cells.Item(31, 9).Value = 2332.7693  # I31: was 2332.8462
$ws.Cells.Item(31, 10).Value = 7245.6665  # J31: was 7932.75
$ws.Cells.Item(31, 11).Value = 2332.7693  # K31: was 2332.8462
$ws.Cells.Item(31, 12).Value = 7245.6665  # L31: was 7932.75
$ws.Cells.Item(31, 13).Value = -2037.7693  # M31: was -2037.8462
$ws.Cells.Item(31, 14).Value = -7835.6665  # N31: was -8522.75
$ws.Cells.Item(34, 8).Value = 4342.591  # H34: was 4466.143
$ws.Cells.Item(34, 9).Value = 2332.7693  # I34: was 2332.8462
$ws.Cells.Item(34, 10).Value = 7245.6665  # J34: was 7932.75
$ws.Cells.Item(34, 11).Value = 2332.7693  # K34: was 2332.8462
$ws.Cells.Item(34, 12).Value = 7245.6665  # L34: was 7932.75
$ws.Cells.Item(34, 13).Value = -2130.7693  # M34: was -2130.8462
$ws.Cells.Item(34, 14).Value = -7649.6665  # N34: was -8336.75
$ws.Cells.Item(54, 8).Value = 34382  # H54: was 34447.2
$ws.Cells.Item(54, 10).Value = 34230.5  # J54: was 34312
$ws.Cells.Item(54, 12).Value = 34230.5  # L54: was 34312
$ws.Cells.Item(54, 14).Value = -35546.5  # N54: was -35628
$ws.Cells.Item(58, 8).Value = 3666.6667  # H58: was 2299.3333
$ws.Cells.Item(58, 9).Value = 1000  # I58: was 949
$ws.Cells.Item(58, 11).Value = 1000  # K58: was 949
$ws.Cells.Item(58, 13).Value = -797  # M58: was -746
$ws.Cells.Item(82, 8).Value = 22333.334  # H82: was 21750
$ws.Cells.Item(82, 10).Value = 22333.334  # J82: was 21750
$ws.Cells.Item(82, 12).Value = 22333.334  # L82: was 21750
$ws.Cells.Item(82, 14).Value = -23055.334  # N82: was -22472
$ws.Cells.Item(85, 8).Value = 22333.334  # H85: was 21750
$ws.Cells.Item(85, 10).Value = 22333.334  # J85: was 21750
$ws.Cells.Item(85, 12).Value = 22333.334  # L85: was 21750
$ws.Cells.Item(85, 14).Value = -24829.334  # N85: was -24246
$ws.Cells.Item(107, 8).Value = 4673.75  # H107: was 4631.6665
$ws.Cells.Item(107, 10).Value = 4673.75  # J107: was 4631.6665
$ws.Cells.Item(107, 12).Value = 4673.75  # L107: was 4631.6665
$ws.Cells.Item(107, 14).Value = -8513.75  # N107: was -8471.666499999999
$ws.Cells.Item(132, 8).Value = 288529  # H132: was 504177.5
$ws.Cells.Item(132, 9).Value = 667270.3  # I132: was 2000012
$ws.Cells.Item(132, 10).Value = 4473  # J132: was 5566
$ws.Cells.Item(132, 11).Value = 2001810.9  # K132: was 6000036
$ws.Cells.Item(132, 12).Value = 13419  # L132: was 16698
$ws.Cells.Item(132, 13).Value = -1999280.9  # M132: was -5997506
$ws.Cells.Item(132, 14).Value = -18479  # N132: was -21758
$ws.Cells.Item(134, 8).Value = 1916.5  # H134: was 1963.5454
$ws.Cells.Item(134, 9).Value = 1777.6666  # I134: was 1825
$ws.Cells.Item(134, 11).Value = 5332.9998  # K134: was 5475
$ws.Cells.Item(134, 13).Value = -2797.9998  # M134: was -2940
$ws.Cells.Item(136, 8).Value = 3666.6667  # H136: was 2299.3333
$ws.Cells.Item(136, 9).Value = 1000  # I136: was 949
$ws.Cells.Item(136, 11).Value = 3000  # K136: was 2847
$ws.Cells.Item(136, 13).Value = -450  # M136: was -297

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(3, 8).Value = 7391.6665  # H3: was 7543.75
$ws.Cells.Item(76, 8).Value = 0  # H76: was 6666
$ws.Cells.Item(76, 10).Value = 0  # J76: was 6666
$ws.Cells.Item(76, 12).Value = 0  # L76: was 19998
$ws.Cells.Item(76, 14).ClearContents()  # N76: was -20764
$ws.Cells.Item(79, 8).Value = 0  # H79: was 6666
$ws.Cells.Item(79, 10).Value = 0  # J79: was 6666
$ws.Cells.Item(79, 12).Value = 0  # L79: was 19998
$ws.Cells.Item(79, 14).ClearContents()  # N79: was -22650
$ws.Cells.Item(140, 8).Value = 1252.2693  # H140: was 931.381
$ws.Cells.Item(140, 9).Value = 674.65  # I140: was 735.7222
$ws.Cells.Item(140, 10).Value = 3177.6667  # J140: was 2105.3333
$ws.Cells.Item(140, 11).Value = 2023.95  # K140: was 2207.1666
$ws.Cells.Item(140, 12).Value = 9533.000100000001  # L140: was 6315.999899999999
$ws.Cells.Item(140, 13).Value = 3156.05  # M140: was 2972.8334
$ws.Cells.Item(140, 14).Value = -19893.0001  # N140: was -16675.9999

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(70, 8).Value = 77440.664  # H70: was 86779.414
$ws.Cells.Item(70, 9).Value = 90317.48  # I70: was 98681.336
$ws.Cells.Item(70, 10).Value = 3399  # J70: was 3466
$ws.Cells.Item(70, 11).Value = 90317.48  # K70: was 98681.336
$ws.Cells.Item(70, 12).Value = 3399  # L70: was 3466
$ws.Cells.Item(70, 13).Value = -90047.48  # M70: was -98411.336
$ws.Cells.Item(70, 14).Value = -3939  # N70: was -4006
$ws.Cells.Item(73, 8).Value = 77440.664  # H73: was 86779.414
$ws.Cells.Item(73, 9).Value = 90317.48  # I73: was 98681.336
$ws.Cells.Item(73, 10).Value = 3399  # J73: was 3466
$ws.Cells.Item(73, 11).Value = 90317.48  # K73: was 98681.336
$ws.Cells.Item(73, 12).Value = 3399  # L73: was 3466
$ws.Cells.Item(73, 13).Value = -89381.48  # M73: was -97745.336
$ws.Cells.Item(73, 14).Value = -5271  # N73: was -5338
$ws.Cells.Item(100, 8).Value = 94999.5  # H100: was 107499.5
$ws.Cells.Item(100, 10).Value = 94999.5  # J100: was 107499.5
$ws.Cells.Item(100, 12).Value = 94999.5  # L100: was 107499.5
$ws.Cells.Item(100, 14).Value = -97163.5  # N100: was -109663.5
$ws.Cells.Item(102, 8).Value = 2255.5  # H102: was 2498
$ws.Cells.Item(102, 9).Value = 1277.15  # I102: was 1405.0625
$ws.Cells.Item(102, 10).Value = 4212.2  # J102: was 4246.7
$ws.Cells.Item(102, 11).Value = 1277.15  # K102: was 1405.0625
$ws.Cells.Item(102, 12).Value = 4212.2  # L102: was 4246.7
$ws.Cells.Item(102, 13).Value = 344.8499999999999  # M102: was 216.9375
$ws.Cells.Item(102, 14).Value = -7456.2  # N102: was -7490.7
$ws.Cells.Item(113, 8).Value = 3774.5227  # H113: was 3864.262
$ws.Cells.Item(113, 9).Value = 3482.0417  # I113: was 3551.261
$ws.Cells.Item(113, 10).Value = 4125.5  # J113: was 4243.1577
$ws.Cells.Item(113, 11).Value = 3482.0417  # K113: was 3551.261
$ws.Cells.Item(113, 12).Value = 4125.5  # L113: was 4243.1577
$ws.Cells.Item(113, 13).Value = -1312.0417  # M113: was -1381.261
$ws.Cells.Item(113, 14).Value = -8465.5  # N113: was -8583.1577
$ws.Cells.Item(126, 8).Value = 4621.7144  # H126: was 4451.778
$ws.Cells.Item(126, 10).Value = 4890  # J126: was 4631.75
$ws.Cells.Item(126, 12).Value = 14670  # L126: was 13895.25
$ws.Cells.Item(126, 14).Value = -19610  # N126: was -18835.25
$ws.Cells.Item(132, 8).Value = 1699.25  # H132: was 2082.5
$ws.Cells.Item(132, 9).Value = 1539.6  # I132: was 1900
$ws.Cells.Item(132, 10).Value = 2497.5  # J132: was 2995
$ws.Cells.Item(132, 11).Value = 4618.799999999999  # K132: was 5700
$ws.Cells.Item(132, 12).Value = 7492.5  # L132: was 8985
$ws.Cells.Item(132, 13).Value = -2088.799999999999  # M132: was -3170
$ws.Cells.Item(132, 14).Value = -12552.5  # N132: was -14045

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(16, 8).Value = 812.9  # H16: was 926.2632
$ws.Cells.Item(16, 9).Value = 490.5  # I16: was 699.5
$ws.Cells.Item(16, 10).Value = 1296.5  # J16: was 1315
$ws.Cells.Item(16, 11).Value = 490.5  # K16: was 699.5
$ws.Cells.Item(16, 12).Value = 1296.5  # L16: was 1315
$ws.Cells.Item(16, 13).Value = -320.5  # M16: was -529.5
$ws.Cells.Item(16, 14).Value = -1636.5  # N16: was -1655
$ws.Cells.Item(54, 8).Value = 33886.89  # H54: was 38331.332
$ws.Cells.Item(54, 9).Value = 29997.5  # I54: was 39995
$ws.Cells.Item(54, 10).Value = 34998.145  # J54: was 37499.5
$ws.Cells.Item(54, 11).Value = 29997.5  # K54: was 39995
$ws.Cells.Item(54, 12).Value = 34998.145  # L54: was 37499.5
$ws.Cells.Item(54, 13).Value = -29353.5  # M54: was -39351
$ws.Cells.Item(54, 14).Value = -36286.145  # N54: was -38787.5
$ws.Cells.Item(58, 8).Value = 25000  # H58: was 0
$ws.Cells.Item(58, 10).Value = 25000  # J58: was 0
$ws.Cells.Item(58, 12).Value = 25000  # L58: was 0
$ws.Cells.Item(58, 14).Value = -25520  # N58: was None
$ws.Cells.Item(61, 8).Value = 2865.1875  # H61: was 3088.7856
$ws.Cells.Item(61, 9).Value = 2278.2856  # I61: was 2669.6
$ws.Cells.Item(61, 11).Value = 2278.2856  # K61: was 2669.6
$ws.Cells.Item(61, 13).Value = -2076.2856  # M61: was -2467.6
$ws.Cells.Item(113, 8).Value = 2865.1875  # H113: was 3088.7856
$ws.Cells.Item(113, 9).Value = 2278.2856  # I113: was 2669.6
$ws.Cells.Item(113, 11).Value = 2278.2856  # K113: was 2669.6
$ws.Cells.Item(113, 13).Value = -108.2856000000002  # M113: was -499.5999999999999
$ws.Cells.Item(132, 8).Value = 12963.587  # H132: was 13313.845
$ws.Cells.Item(132, 9).Value = 10883.308  # I132: was 11698.583
$ws.Cells.Item(132, 10).Value = 15667.95  # J132: was 15159.857
$ws.Cells.Item(132, 11).Value = 32649.924  # K132: was 35095.749
$ws.Cells.Item(132, 12).Value = 47003.85000000001  # L132: was 45479.571
$ws.Cells.Item(132, 13).Value = -30119.924  # M132: was -32565.749
$ws.Cells.Item(132, 14).Value = -52063.85000000001  # N132: was -50539.571

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(122, 8).Value = 999999  # H122: was 253247
$ws.Cells.Item(122, 9).Value = 999999  # I122: was 501999.5
$ws.Cells.Item(122, 10).Value = 0  # J122: was 4494.5
$ws.Cells.Item(122, 11).Value = 2999997  # K122: was 1505998.5
$ws.Cells.Item(122, 12).Value = 0  # L122: was 13483.5
$ws.Cells.Item(122, 13).Value = -2997547  # M122: was -1503548.5
$ws.Cells.Item(122, 14).ClearContents()  # N122: was -18383.5
